# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff", and the associated timestamps were refreshed to the
# moment the handoff package was produced. The same status/timestamp
# fields are mirrored on the Overview summary sheet and on each
# per-locale detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" --------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ---------------------------------------
$overview.Range("G2").Value = "2016-09-01 12:43:41"
$zhcn.Range("H2").Value     = "2016-09-01 12:43:37"
$dede.Range("H2").Value     = "2016-09-01 12:43:41"

# --- The wider "Ready for handoff" text bumps the Status column's
#     auto-fit width on every sheet that shows it. ------------------------
$overview.Range("E1").EntireColumn.ColumnWidth = 16.33
$overview.Range("F1").EntireColumn.ColumnWidth = 16.33
$zhcn.Range("C1").EntireColumn.ColumnWidth     = 16.33
$dede.Range("C1").EntireColumn.ColumnWidth     = 16.33
